$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-08 Tuesday" "2024-10-09 Wednesday"

Replace-Text "96×47=4512" "72×15=1080"
Replace-Text "48×29=1392" "32×25=800"
Replace-Text "76×41=3116" "18×27=486"
Replace-Text "73×59=4307" "86×14=1204"
Replace-Text "18×76=1368" "32×35=1120"

Replace-Text "94×67=6298" "45×56=2520"
Replace-Text "67×23=1541" "81×29=2349"
Replace-Text "67×16=1072" "93×22=2046"
Replace-Text "98×28=2744" "78×86=6708"
Replace-Text "17×41=697" "63×97=6111"

Replace-Text "26×82=2132" "73×22=1606"
Replace-Text "70×90=6300" "41×66=2706"
Replace-Text "17×23=391" "44×42=1848"
Replace-Text "44×95=4180" "83×26=2158"
Replace-Text "94×94=8836" "97×85=8245"

Replace-Text "72×13=936" "62×56=3472"
Replace-Text "94×59=5546" "91×14=1274"
Replace-Text "68×63=4284" "47×54=2538"
Replace-Text "23×55=1265" "75×33=2475"
Replace-Text "48×78=3744" "92×38=3496"

Replace-Text "16×60=960" "29×87=2523"
Replace-Text "38×58=2204" "63×31=1953"
Replace-Text "85×85=7225" "69×57=3933"
Replace-Text "47×93=4371" "83×31=2573"
Replace-Text "17×77=1309" "97×82=7954"
